$wb = $excel.ActiveWorkbook

# Sheet "zh-cn" row 3: Correspond Handoff Datetime (E3) and Correspond Handback DateTime (H3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 12:36:27"
$wsZhCn.Range("H3").Value = "2016-03-20 12:36:48"

# Sheet "de-de" row 3: Correspond Handoff Datetime (E3) and Correspond Handback DateTime (H3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 12:36:30"
$wsDeDe.Range("H3").Value = "2016-03-20 12:36:53"
